# Cryptos sheet refresh: update Price (D) / Volume(1h) (E) text cells.
# Numeric-looking Price values are written as "'" + value so Excel keeps
# them as Text (matching the source's inlineStr cells) instead of silently
# coercing them to numbers (which would e.g. turn "15.60" into 15.6).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '27.885.48'
$ws.Range("E2").Value = '  -0.14%  '

# Row 3
$ws.Range("D3").Value = '1.634.19'
$ws.Range("E3").Value = '  -0.81%  '

# Row 4
$ws.Range("E4").Value = '  +0.01%  '

# Row 5
$ws.Range("D5").Value = "'" + '212.07'
$ws.Range("E5").Value = '  -0.80%  '

# Row 6
$ws.Range("E6").Value = '  -0.77%  '

# Row 7
$ws.Range("E7").Value = '  +0.02%  '

# Row 8
$ws.Range("D8").Value = "'" + '23.19'
$ws.Range("E8").Value = '  -1.54%  '

# Row 9
$ws.Range("E9").Value = '  -3.09%  '

# Row 10
$ws.Range("E10").Value = '  -0.35%  '

# Row 11
$ws.Range("D11").Value = "'" + '0.0881'
$ws.Range("E11").Value = '  +0.98%  '

# Row 12
$ws.Range("D12").Value = '1.866.31'
$ws.Range("E12").Value = '  -0.74%  '

# Row 13
$ws.Range("D13").Value = '1.632.36'
$ws.Range("E13").Value = '  -0.85%  '

# Row 14
$ws.Range("E14").Value = '  -0.57%  '

# Row 15
$ws.Range("D15").Value = "'" + '0.567'
$ws.Range("E15").Value = '  +0.62%  '

# Row 16
$ws.Range("D16").Value = "'" + '65.13'
$ws.Range("E16").Value = '  -0.78%  '

# Row 17
$ws.Range("D17").Value = '27.894.28'
$ws.Range("E17").Value = '  -0.10%  '

# Row 18
$ws.Range("D18").Value = "'" + '229.41'
$ws.Range("E18").Value = '  -1.28%  '

# Row 19
$ws.Range("D19").Value = '0.0₃0719'
$ws.Range("E19").Value = '  -0.56%  '

# Row 20
$ws.Range("D20").Value = "'" + '7.49'
$ws.Range("E20").Value = '  -2.08%  '

# Row 21
$ws.Range("D21").Value = "'" + '0.999'
$ws.Range("E21").Value = '  -0.10%  '

# Row 22
$ws.Range("D22").Value = "'" + '4.36'
$ws.Range("E22").Value = '  -0.74%  '

# Row 23
$ws.Range("D23").Value = "'" + '10.31'
$ws.Range("E23").Value = '  -3.93%  '

# Row 24
$ws.Range("E24").Value = '  -3.59%  '

# Row 25
$ws.Range("D25").Value = "'" + '152.92'
$ws.Range("E25").Value = '  +0.46%  '

# Row 26
$ws.Range("D26").Value = "'" + '6.95'
$ws.Range("E26").Value = '  +0.71%  '

# Row 27
$ws.Range("D27").Value = "'" + '15.60'
$ws.Range("E27").Value = '  -0.81%  '

# Row 28
$ws.Range("E28").Value = '  -0.72%  '

# Row 29
$ws.Range("E29").Value = '  +0.01%  '

# Row 30
$ws.Range("E30").Value = '  -1.00%  '

# Row 31
$ws.Range("E31").Value = '  -0.76%  '

# Row 32
$ws.Range("E32").Value = '  +0.73%  '

# Row 33
$ws.Range("D33").Value = '1.404.89'
$ws.Range("E33").Value = '  -3.28%  '

# Row 34
$ws.Range("D34").Value = "'" + '3.06'
$ws.Range("E34").Value = '  -1.82%  '

# Row 35
$ws.Range("D35").Value = "'" + '1.56'
$ws.Range("E35").Value = '  +0.70%  '

# Row 36
$ws.Range("D36").Value = "'" + '0.999'
$ws.Range("E36").Value = '  +8.45%  '

# Row 37
$ws.Range("E37").Value = '  +1.55%  '

# Row 38
$ws.Range("E38").Value = '  +0.44%  '

# Row 39
$ws.Range("D39").Value = "'" + '0.560'
$ws.Range("E39").Value = '  -0.10%  '

# Row 40
$ws.Range("D40").Value = "'" + '0.870'
$ws.Range("E40").Value = '  -2.30%  '

# Row 41
$ws.Range("E41").Value = '  +0.07%  '

# Row 42
$ws.Range("E42").Value = '  -0.03%  '

# Row 43
$ws.Range("D43").Value = "'" + '66.82'
$ws.Range("E43").Value = '  -3.66%  '

# Row 44
$ws.Range("D44").Value = "'" + '5.49'
$ws.Range("E44").Value = '  +2.16%  '

# Row 45
$ws.Range("D45").Value = "'" + '1.80'
$ws.Range("E45").Value = '  +1.09%  '

# Row 46
$ws.Range("E46").Value = '  -1.57%  '

# Row 47
$ws.Range("D47").Value = '1.775.91'
$ws.Range("E47").Value = '  -0.78%  '

# Row 48
$ws.Range("D48").Value = "'" + '87.63'
$ws.Range("E48").Value = '  -1.45%  '

# Row 49
$ws.Range("E49").Value = '  -0.72%  '

# Row 50
$ws.Range("E50").Value = '  -0.12%  '

# Row 51
$ws.Range("D51").Value = "'" + '7.54'
$ws.Range("E51").Value = '  -2.19%  '

